# fixes numbering of chapter 3
#
# The deck's "x / 8" slide-number footers were left over from an 8-slide
# version of the chapter; the chapter now only has 6 slides, so every
# footer needs to read "x / 6" instead. Slide 3 (the section divider) has
# no slide-number placeholder, so only slides 1, 2, 4, 5 and 6 need fixing.

$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

# Slides 1, 2 and 4: "N / 8" -> "N / " + "6" (split into two runs, the
# second run carrying the corrected total page count).
$fixups = @{ 1 = "1"; 2 = "2"; 4 = "4" }

foreach ($slideIndex in $fixups.Keys) {
    $slide = $p.Slides.Item($slideIndex)
    $shp = Get-ShapeByName $slide "Slide Number Placeholder 4"
    $tr = $shp.TextFrame.TextRange
    $n = $fixups[$slideIndex]
    $tr.Text = "$n / "
    $tr.InsertAfter("6") | Out-Null
}

# Slide 5: "7 / 8" -> "5 / 6" (the slide itself was renumbered from 7 to 5),
# split across four runs: "5", " ", "/ ", "6".
$slide5 = $p.Slides.Item(5)
$shp5 = Get-ShapeByName $slide5 "Slide Number Placeholder 4"
$tr5 = $shp5.TextFrame.TextRange
$tr5.Text = "5"
$tr5.InsertAfter(" ") | Out-Null
$tr5.InsertAfter("/ ") | Out-Null
$tr5.InsertAfter("6") | Out-Null

# Slide 6: "6 / " + "6" (two runs) -> a single run "6 / 6". The
# concatenated text already reads "6 / 6", so assigning that directly
# would be treated as a no-op (and an intermediate value that shares a
# prefix/suffix with it would only patch around the existing run split
# instead of collapsing it) - route through an unrelated placeholder
# value first so the two runs actually get merged into one.
$slide6 = $p.Slides.Item(6)
$shp6 = Get-ShapeByName $slide6 "Slide Number Placeholder 4"
$tr6 = $shp6.TextFrame.TextRange
$tr6.Text = "placeholder"
$tr6.Text = "6 / 6"
